# Applies the two changes from the commit:
#  1. The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") gets a new
#     built-in table style applied (tableStyleId swaps from the custom
#     "Table_0" style {A380631B-...} to the built-in style {B125649A-...}).
#  2. The presentation's theme (ppt/theme/theme2.xml, the theme actually
#     used by the slide master / whole deck) is switched from the
#     "Integral" / "Red Violet" colour scheme back to the plain default
#     "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 ---------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B125649A-5BF7-40A1-AB4C-79D40CEFFC8E}")
    }
}

# --- 2) Restore the "Office Theme" colour scheme -----------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

try { $theme.Name = "Office Theme" } catch { }
try { $colors.Name = "Office" } catch { }

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
